{"js": "// Replace the date line and each \"A\u00d7B=C\" multiplication-table answer\n// with the new values from the commit, preserving all run/paragraph\n// formatting by doing an in-place text replacement on the matched range.\nconst replacements = [\n  [\"2025-02-08 Saturday\", \"2025-02-09 Sunday\"],\n  [\"977\u00d73=2931\", \"141\u00d75=705\"],\n  [\"689\u00d72=1378\", \"629\u00d74=2516\"],\n  [\"651\u00d79=5859\", \"988\u00d73=2964\"],\n  [\"806\u00d77=5642\", \"813\u00d79=7317\"],\n  [\"696\u00d77=4872\", \"479\u00d73=1437\"],\n  [\"497\u00d76=2982\", \"124\u00d78=992\"],\n  [\"305\u00d75=1525\", \"265\u00d78=2120\"],\n  [\"671\u00d74=2684\", \"525\u00d76=3150\"],\n  [\"266\u00d72=532\", \"219\u00d75=1095\"],\n  [\"169\u00d75=845\", \"435\u00d79=3915\"],\n  [\"606\u00d74=2424\", \"823\u00d73=2469\"],\n  [\"593\u00d77=4151\", \"477\u00d79=4293\"],\n  [\"844\u00d72=1688\", \"892\u00d79=8028\"],\n  [\"712\u00d75=3560\", \"124\u00d72=248\"],\n  [\"290\u00d75=1450\", \"514\u00d79=4626\"],\n  [\"974\u00d72=1948\", \"559\u00d76=3354\"],\n  [\"121\u00d73=363\", \"540\u00d79=4860\"],\n  [\"189\u00d72=378\", \"666\u00d79=5994\"],\n  [\"713\u00d76=4278\", \"401\u00d72=802\"],\n  [\"605\u00d76=3630\", \"150\u00d74=600\"],\n  [\"646\u00d77=4522\", \"566\u00d77=3962\"],\n  [\"305\u00d79=2745\", \"944\u00d75=4720\"],\n  [\"435\u00d76=2610\", \"710\u00d79=6390\"],\n  [\"773\u00d72=1546\", \"105\u00d75=525\"],\n  [\"806\u00d75=4030\", \"607\u00d77=4249\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"A\u00d7B=C\" multiplication-table answer\n# with the new values from the commit, preserving all run/paragraph\n# formatting by using Find/Replace over the whole document content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-08 Saturday\", \"2025-02-09 Sunday\"),\n    @(\"977\u00d73=2931\", \"141\u00d75=705\"),\n    @(\"689\u00d72=1378\", \"629\u00d74=2516\"),\n    @(\"651\u00d79=5859\", \"988\u00d73=2964\"),\n    @(\"806\u00d77=5642\", \"813\u00d79=7317\"),\n    @(\"696\u00d77=4872\", \"479\u00d73=1437\"),\n    @(\"497\u00d76=2982\", \"124\u00d78=992\"),\n    @(\"305\u00d75=1525\", \"265\u00d78=2120\"),\n    @(\"671\u00d74=2684\", \"525\u00d76=3150\"),\n    @(\"266\u00d72=532\", \"219\u00d75=1095\"),\n    @(\"169\u00d75=845\", \"435\u00d79=3915\"),\n    @(\"606\u00d74=2424\", \"823\u00d73=2469\"),\n    @(\"593\u00d77=4151\", \"477\u00d79=4293\"),\n    @(\"844\u00d72=1688\", \"892\u00d79=8028\"),\n    @(\"712\u00d75=3560\", \"124\u00d72=248\"),\n    @(\"290\u00d75=1450\", \"514\u00d79=4626\"),\n    @(\"974\u00d72=1948\", \"559\u00d76=3354\"),\n    @(\"121\u00d73=363\", \"540\u00d79=4860\"),\n    @(\"189\u00d72=378\", \"666\u00d79=5994\"),\n    @(\"713\u00d76=4278\", \"401\u00d72=802\"),\n    @(\"605\u00d76=3630\", \"150\u00d74=600\"),\n    @(\"646\u00d77=4522\", \"566\u00d77=3962\"),\n    @(\"305\u00d79=2745\", \"944\u00d75=4720\"),\n    @(\"435\u00d76=2610\", \"710\u00d79=6390\"),\n    @(\"773\u00d72=1546\", \"105\u00d75=525\"),\n    @(\"806\u00d75=4030\", \"607\u00d77=4249\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Replacement.Text = $newText\n    $rng.Find.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
